$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 110
$ws.Range("I2").Value = 276
$ws.Range("J2").Value = 1266
$ws.Range("K2").Value = 9
$ws.Range("L2").Value = 338
$ws.Range("M2").Value = 20
$ws.Range("N2").Value = 214
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 2
$ws.Range("Q2").Value = 2
$ws.Range("R2").Value = 19
$ws.Range("S2").Value = 123
$ws.Range("T2").Value = 191
$ws.Range("U2").Value = 13
$ws.Range("V2").Value = 1886
$ws.Range("W2").Value = 1
$ws.Range("X2").Value = 1927
$ws.Range("Y2").Value = 2
$ws.Range("Z2").Value = 32
$ws.Range("AA2").Value = 17
